$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values are purely numeric-looking (e.g. "0.682")
# must be forced to Text format first, otherwise Excel auto-converts them
# to floating-point numbers and mangles the formatting (trailing zeros, etc).
# We reset Style back to "Normal" afterwards so no stray number-format / style
# is left behind on the cell (matching the original workbook, which has none).

$ws.Range('D2').Value = '44.253.69'
$ws.Range('E2').Value = '  +1.91%  '
$ws.Range('D3').Value = '2.362.90'
$ws.Range('E3').Value = '  -0.43%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('B5').Value = 'XRP'
$ws.Range('C5').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.682'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.78%  '
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '243.92'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.56%  '
$ws.Range('E7').Value = '  +4.45%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +24.76%  '
$ws.Range('E10').Value = '  +5.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '32.23'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +21.55%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.49'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +19.78%  '
$ws.Range('E13').Value = '  +2.18%  '
$ws.Range('D14').Value = '2.713.21'
$ws.Range('E14').Value = '  -0.46%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '16.95'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +7.43%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.917'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +7.26%  '
$ws.Range('D17').Value = '2.361.96'
$ws.Range('E17').Value = '  -0.58%  '
$ws.Range('D18').Value = '44.211.73'
$ws.Range('E18').Value = '  +1.79%  '
$ws.Range('E19').Value = '  +4.44%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.81'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.99%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '78.48'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '256.51'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.95%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.59'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.37%  '
$ws.Range('E25').Value = '  -5.31%  '
$ws.Range('E26').Value = '  +7.48%  '
$ws.Range('E27').Value = '  +3.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '22.61'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.24%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '175.20'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.52%  '
$ws.Range('E30').Value = '  +3.78%  '
$ws.Range('E31').Value = '  +3.82%  '
$ws.Range('E32').Value = '  +4.90%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.44'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +8.79%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0757'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +9.24%  '
$ws.Range('E35').Value = '  +5.52%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.90'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.82%  '
$ws.Range('E37').Value = '  +0.66%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.60'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0275'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +7.99%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '19.46'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '9.02'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.71%  '
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('E43').Value = '  +14.81%  '
$ws.Range('E44').Value = '  +3.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.51'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +11.71%  '
$ws.Range('E46').Value = '  +5.16%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '101.60'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.92%  '
$ws.Range('E48').Value = '  -0.23%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.49'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.35%  '
$ws.Range('B50').Value = 'TerraClassic'
$ws.Range('C50').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000209'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.10%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '1.455.78'
$ws.Range('E51').Value = '  +0.13%  '
